$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to stay text so numeric-looking values
# like "112.26" are not auto-converted to numbers by Excel.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "42.518.86"
$ws.Range("E2").Value = "  -0.67%  "
$ws.Range("D3").Value = "2.229.42"
$ws.Range("E3").Value = "  -0.50%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").Value = "112.26"
$ws.Range("E5").Value = "  -1.33%  "
$ws.Range("D6").Value = "293.15"
$ws.Range("E6").Value = "  +8.43%  "
$ws.Range("D7").Value = "0.625"
$ws.Range("E7").Value = "  -0.33%  "
$ws.Range("E8").Value = "  -0.30%  "
$ws.Range("D9").Value = "0.599"
$ws.Range("E9").Value = "  -1.02%  "
$ws.Range("D10").Value = "43.47"
$ws.Range("E10").Value = "  -6.08%  "
$ws.Range("D11").Value = "0.0922"
$ws.Range("E11").Value = "  -0.94%  "
$ws.Range("D12").Value = "54.15"
$ws.Range("E12").Value = "  +0.29%  "
$ws.Range("D13").Value = "8.73"
$ws.Range("E13").Value = "  -4.32%  "
$ws.Range("E14").Value = "  +21.05%  "
$ws.Range("E15").Value = "  -1.30%  "
$ws.Range("D16").Value = "14.96"
$ws.Range("E16").Value = "  -2.48%  "
$ws.Range("D17").Value = "2.564.85"
$ws.Range("E17").Value = "  -0.53%  "
$ws.Range("D18").Value = "2.243.62"
$ws.Range("E18").Value = "  +0.23%  "
$ws.Range("D19").Value = "42.503.28"
$ws.Range("E19").Value = "  -1.07%  "
$ws.Range("D20").Value = "7.22"
$ws.Range("E20").Value = "  +6.81%  "
$ws.Range("E21").Value = "  -1.71%  "
$ws.Range("E22").Value = "  +2.17%  "
$ws.Range("D23").Value = "3.35"
$ws.Range("E23").Value = "  +13.90%  "
$ws.Range("E24").Value = "  +0.87%  "
$ws.Range("D25").Value = "241.06"
$ws.Range("E25").Value = "  +4.22%  "
$ws.Range("D26").Value = "8.85"
$ws.Range("E26").Value = "  -4.65%  "
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  -1.24%  "
$ws.Range("D28").Value = "11.45"
$ws.Range("E28").Value = "  -6.05%  "
$ws.Range("E29").Value = "  -1.85%  "
$ws.Range("D30").Value = "175.41"
$ws.Range("E30").Value = "  +0.96%  "
$ws.Range("D31").Value = "37.00"
$ws.Range("E31").Value = "  -8.23%  "
$ws.Range("D32").Value = "21.71"
$ws.Range("E32").Value = "  +2.70%  "
$ws.Range("E33").Value = "  -4.67%  "
$ws.Range("D34").Value = "0.0879"
$ws.Range("E34").Value = "  -2.63%  "
$ws.Range("D35").Value = "5.67"
$ws.Range("E35").Value = "  +1.73%  "
$ws.Range("E36").Value = "  +3.99%  "
$ws.Range("E37").Value = "  -0.55%  "
$ws.Range("D38").Value = "4.17"
$ws.Range("E38").Value = "  -2.83%  "
$ws.Range("D39").Value = "0.0372"
$ws.Range("E39").Value = "  -0.66%  "
$ws.Range("E40").Value = "  -2.62%  "
$ws.Range("E41").Value = "  -6.84%  "
$ws.Range("D42").Value = "71.13"
$ws.Range("E42").Value = "  +0.35%  "
$ws.Range("E43").Value = "  -1.89%  "
$ws.Range("E44").Value = "  -0.04%  "
$ws.Range("D45").Value = "12.29"
$ws.Range("E45").Value = "  -6.99%  "
$ws.Range("E46").Value = "  -2.61%  "
$ws.Range("D47").Value = "5.41"
$ws.Range("E47").Value = "  -4.52%  "
$ws.Range("E48").Value = "  +2.24%  "
$ws.Range("D49").Value = "8.51"
$ws.Range("E49").Value = "  +0.79%  "
$ws.Range("D50").Value = "102.09"
$ws.Range("E50").Value = "  +1.38%  "
$ws.Range("E51").Value = "  -1.45%  "

# Restore the original (default) cell style on column D now that the
# text values are safely stored, matching the unstyled source cells.
$priceRange.ClearFormats()
